# Regenerate sval data to filter save games: update computed columns B-E and G (sum)
# for rows 2-13 on the active worksheet. Column F (Win) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 1.445647641019636;     C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 4.327115817150455 }
    3  = @{ B = 3.272327238179451;     C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987;  G = 8.656069925401464 }
    4  = @{ B = 0.003078177322033415;  C = 0.3048912486333797; D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.562449902544138 }
    5  = @{ B = 3.272327238179451;     C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    6  = @{ B = 3.272327238179451;     C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987;  G = 8.656069925401464 }
    7  = @{ B = 0.1169995834814548;    C = 0.3048912486333797; D = 0.1496068669990043; E = 0.5333859586016987;  G = 1.104883657715537 }
    8  = @{ B = 0.1169995834814548;    C = 0.3048912486333797; D = 0.7210945179870265; E = 13.86384647080068;   G = 15.00683182090255 }
    9  = @{ B = 0.000001174341637932841; C = 0.04103571897497393; D = 18.71679738969934;  E = 14773364.14517103;   G = 14773382.90300531 }
    10 = @{ B = 3.272327238179451;     C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    11 = @{ B = 3.272327238179451;     C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    12 = @{ B = 1.445647641019636;     C = 1.626987699542094;  D = 3.223369029078222;  E = 13.86384647080068;   G = 20.15985084044064 }
    13 = @{ B = 3.272327238179451;     C = 2919.202174992006;  D = 19575605.8673771;   E = 2797.565817734744;   G = 19581325.90769707 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B   # B = TB
    $ws.Cells.Item($row, 3).Value = $vals.C   # C = d2S
    $ws.Cells.Item($row, 4).Value = $vals.D   # D = K
    $ws.Cells.Item($row, 5).Value = $vals.E   # E = IP
    $ws.Cells.Item($row, 7).Value = $vals.G   # G = sum
}
